$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.904.85"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.578.36"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  -1.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "190.50"
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("E7").Value = "  -1.99%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.576.71"
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("E10").Value = "  +2.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.663"
$ws.Range("E11").Value = "  +0.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.81"
$ws.Range("E12").Value = "  -3.60%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000306"
$ws.Range("E13").Value = "  +5.06%  "
$ws.Range("E14").Value = "  -0.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.156.14"
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.94"
$ws.Range("E16").Value = "  +3.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.577.89"
$ws.Range("E17").Value = "  -0.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.895.16"
$ws.Range("E18").Value = "  -0.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.68"
$ws.Range("E19").Value = "  +1.55%  "
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("E21").Value = "  -0.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "478.69"
$ws.Range("E22").Value = "  -3.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "19.36"
$ws.Range("E23").Value = "  +10.91%  "
$ws.Range("E24").Value = "  -5.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.16"
$ws.Range("E25").Value = "  +6.26%  "
$ws.Range("E26").Value = "  -1.66%  "
$ws.Range("E27").Value = "  -3.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.06"
$ws.Range("E28").Value = "  -0.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.42"
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.37"
$ws.Range("E30").Value = "  +0.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.68"
$ws.Range("E31").Value = "  +2.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.23"
$ws.Range("E32").Value = "  +0.38%  "
$ws.Range("E33").Value = "  +1.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "66.22"
$ws.Range("E34").Value = "  +1.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "582.51"
$ws.Range("E35").Value = "  -6.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "38.99"
$ws.Range("E36").Value = "  +2.73%  "
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0802"
$ws.Range("E38").Value = "  -1.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.396"
$ws.Range("E39").Value = "  -2.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.24"
$ws.Range("E40").Value = "  +20.40%  "
$ws.Range("E41").Value = "  -5.77%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.47"
$ws.Range("E42").Value = "  -4.67%  "
$ws.Range("E43").Value = "  +7.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.230.21"
$ws.Range("E44").Value = "  -2.62%  "
$ws.Range("E45").Value = "  +0.44%  "
$ws.Range("E46").Value = "  -0.25%  "
$ws.Range("E47").Value = "  +0.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.37"
$ws.Range("E48").Value = "  +3.11%  "
$ws.Range("E49").Value = "  +0.71%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.14"
$ws.Range("E51").Value = "  -5.21%  "
